# Update "艺术表演团体国内演出观众人次" sheet:
# - Remove the 2008年 and 2009年 rows (old rows 2 and 3), shifting all
#   subsequent year rows up by two.
# - Append a new 2021年 row at the end (row 13 after the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two oldest rows (2008, 2009). Deleting row 2 twice removes
# both, shifting everything else up.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# After the deletion, the last data row (2020年) is row 12. Add the new
# 2021年 row right after it, at row 13.
$newRow = 13

# Match the style of the other "A" column year labels (bold/centered/
# bordered) by copying the formatting from the cell directly above.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = "2021年"
$ws.Cells.Item($newRow, 2).Value = 825.2
$ws.Cells.Item($newRow, 3).Value = ""
$ws.Cells.Item($newRow, 4).Value = ""
$ws.Cells.Item($newRow, 5).Value = 757109.9
$ws.Cells.Item($newRow, 6).Value = 742788.1
$ws.Cells.Item($newRow, 7).Value = 858037.2
$ws.Cells.Item($newRow, 8).Value = 154857.8
$ws.Cells.Item($newRow, 9).Value = 46195.4
$ws.Cells.Item($newRow, 10).Value = ""
$ws.Cells.Item($newRow, 11).Value = 191902.3
$ws.Cells.Item($newRow, 12).Value = 736132.7
$ws.Cells.Item($newRow, 13).Value = 185246.8
$ws.Cells.Item($newRow, 14).Value = ""
$ws.Cells.Item($newRow, 15).Value = ""
$ws.Cells.Item($newRow, 16).Value = ""
$ws.Cells.Item($newRow, 17).Value = ""
$ws.Cells.Item($newRow, 18).Value = 22977.2
$ws.Cells.Item($newRow, 19).Value = ""
$ws.Cells.Item($newRow, 20).Value = 928035
$ws.Cells.Item($newRow, 21).Value = ""
$ws.Cells.Item($newRow, 22).Value = 16067.3
